# Insert a new data row before the current row 390 (Perejil / Vega Modelo de Temuco)
# This shifts the existing rows 390-463 down to 391-464, and the new row 390 becomes
# a fresh record (same master-data columns as the former row 390, but with its own
# date and volume figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 390:463 down by inserting a new row at 390
$ws.Rows.Item(390).Insert()

# Populate the newly inserted row 390 with the new record
$ws.Cells.Item(390, 1).Value = 10
$ws.Cells.Item(390, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(390, 3).Value = "La Araucanía"
$ws.Cells.Item(390, 4).Value = 45015
$ws.Cells.Item(390, 5).Value = 9
$ws.Cells.Item(390, 6).Value = 100112044
$ws.Cells.Item(390, 7).Value = "Perejil"
$ws.Cells.Item(390, 8).Value = "Sin especificar"
$ws.Cells.Item(390, 9).Value = "Primera"
$ws.Cells.Item(390, 10).Value = 50
$ws.Cells.Item(390, 11).Value = 4000
$ws.Cells.Item(390, 12).Value = 4000
$ws.Cells.Item(390, 13).Value = 4000
$ws.Cells.Item(390, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(390, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(390, 16).Value = 1333
$ws.Cells.Item(390, 17).Value = 3
$ws.Cells.Item(390, 18).Value = "Hortaliza"
